$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.961.85'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -1.03%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.008.23'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -1.89%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '225.93'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.06%  '
$ws.Range("E6").Value = '  -0.91%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '55.13'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.49%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.372'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -3.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0775'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -4.04%  '
$ws.Range("E11").Value = '  -4.16%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.305.47'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.97%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.02'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -3.15%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.69'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -3.96%  '
$ws.Range("E15").Value = '  -2.46%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.13'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.11%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.008.63'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.58%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '36.924.01'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.88%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.15'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.42%  '
$ws.Range("E20").Value = '  -2.15%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0808'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -4.35%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '222.59'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.33%  '
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("E24").Value = '  +1.85%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.16'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -4.79%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.49'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.38%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.90'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -6.42%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.126'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.51'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.89%  '
$ws.Range("E30").Value = '  -6.27%  '
$ws.Range("E31").Value = '  -1.06%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.39'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -3.11%  '
$ws.Range("E33").Value = '  -1.76%  '
$ws.Range("E34").Value = '  -2.26%  '
$ws.Range("E35").Value = '  -2.89%  '
$ws.Range("E36").Value = '  +2.20%  '
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.15'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.61%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.33'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.35%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.462.86'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.66%  '
$ws.Range("E41").Value = '  -4.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '94.17'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.78%  '
$ws.Range("B43").Value = 'Cronos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0909'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.19%  '
$ws.Range("B44").Value = 'HuobiToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.76'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -4.13%  '
$ws.Range("B45").Value = 'FTXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.19'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +12.17%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.11'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.82%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '15.81'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -5.61%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.994'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.84%  '
$ws.Range("E49").Value = '  -2.10%  '
$ws.Range("E50").Value = '  -1.34%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.195.60'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.87%  '
